$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Update the displayed text of the hyperlink cell G7 (keep the existing
# hyperlink relationship target untouched, only the visible text changes)
$ws.Range("G7").Value = "https://github.com/makersmakingchange/Beverage_Can_Opener/tree/main/Build_Files/3D_Print_Files"

# Update the selected/active cell shown when the workbook is reopened
$ws.Range("G7").Select()
